$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# "Generate Report for Handback" - refresh the localization-status
# report: the handback for the de-de xliff is now in, so every sheet's
# status moves from "Ready for handoff" to "Handed back: in sync with
# en-US", the per-language Handback datetime stamps tick forward, the
# (now resolved) handback-version-mismatch error text is cleared, and
# the Status/Error Detail columns get widened/narrowed to fit the new
# text.
# -----------------------------------------------------------------

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("K2").Value = "2016-09-06 08:01:14"
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8333333333333

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("K2").Value = "2016-09-06 08:01:34"
$wsDeDe.Range("P2").Value = ""
$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8333333333333
